# Add a new run containing two spaces right after the existing
# "Types of Machine Learning: " + " " runs, without disturbing the
# existing runs (so a brand-new <w:r> is created rather than merging
# text into an existing run).

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute(
    "Types of Machine Learning:  ", `
    $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the end of the located text (still before the
    # paragraph mark) and insert a new run of two spaces there.
    $rng.Collapse(0)
    $rng.InsertAfter("  ")
}
